$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("K2").Value = 1.83
$ws.Range("M2").Value = 1.13
$ws.Range("N2").Value = 6
$ws.Range("S2").Value = 2.88
$ws.Range("T2").Value = 1.4
$ws.Range("L3").Value = 2.88
$ws.Range("G5").Value = 2.55
$ws.Range("H5").Value = 3.4
$ws.Range("I5").Value = 2.47
$ws.Range("K5").Value = 2.18
$ws.Range("L5").Value = 3.05
$ws.Range("S5").Value = 1.62
$ws.Range("T5").Value = 2.02
$ws.Range("W5").Value = 2.45
$ws.Range("AA5").Value = 1.52
$ws.Range("AB5").Value = 2.22
$ws.Range("AC5").Value = 11
$ws.Range("AD5").Value = 15
$ws.Range("AF5").Value = 30
$ws.Range("AH5").Value = 23
$ws.Range("AI5").Value = 12.5
$ws.Range("AJ5").Value = 6.8
$ws.Range("AK5").Value = 11.5
$ws.Range("AL5").Value = 40
$ws.Range("AM5").Value = 10.25
$ws.Range("AS5").Value = 250
$ws.Range("G6").Value = 1.67
$ws.Range("H6").Value = 3.8
$ws.Range("I6").Value = 4.33
$ws.Range("L6").Value = 4.5
$ws.Range("M6").Value = 1.04
$ws.Range("N6").Value = 9
$ws.Range("O6").Value = 1.22
$ws.Range("P6").Value = 4
$ws.Range("S6").Value = 1.73
$ws.Range("T6").Value = 2.08
$ws.Range("W6").Value = 2.75
$ws.Range("X6").Value = 1.4
$ws.Range("Y6").Value = 1.33
$ws.Range("AD6").Value = 9
$ws.Range("AJ6").Value = 7.5
$ws.Range("AN6").Value = 23
$ws.Range("AP6").Value = 41
$ws.Range("AR6").Value = 34
$ws.Range("G7").Value = 2.3
$ws.Range("I7").Value = 2.88
$ws.Range("L7").Value = 3.5
$ws.Range("M7").Value = 1.06
$ws.Range("N7").Value = 8
$ws.Range("O7").Value = 1.33
$ws.Range("P7").Value = 3.25
$ws.Range("S7").Value = 2.05
$ws.Range("T7").Value = 1.75
$ws.Range("Y7").Value = 1.44
$ws.Range("Z7").Value = 2.63
$ws.Range("AA7").Value = 1.83
$ws.Range("AB7").Value = 1.83
$ws.Range("AD7").Value = 11
$ws.Range("AF7").Value = 21
$ws.Range("AI7").Value = 9
$ws.Range("AN7").Value = 15
$ws.Range("AO7").Value = 12
$ws.Range("H8").Value = 3.5
$ws.Range("I8").Value = 3.7
$ws.Range("K8").Value = 2.1
$ws.Range("S8").Value = 2.03
$ws.Range("T8").Value = 1.78
$ws.Range("Y8").Value = 1.44
$ws.Range("Z8").Value = 2.63
$ws.Range("AA8").Value = 1.83
$ws.Range("AB8").Value = 1.83
$ws.Range("AC8").Value = 7
$ws.Range("AI8").Value = 9.5
$ws.Range("AM8").Value = 10
$ws.Range("AO8").Value = 13
